# data: merge Palghar, Vasai, Navi Mumbai into Mumbai
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows in column B (CITY) that currently contain "NEW MUMBAI", "VASAI", or "PALGHAR"
# get merged into "MUMBAI".
$rows = @(27, 30, 54, 180, 199, 205, 227, 230, 260)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 2).Value = "MUMBAI"
}
